$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.9887275
$ws.Range("H2").Value = 39.977455
$ws.Range("I2").Value = 0.2607251058521317
$ws.Range("J2").Value = 0.2185266677421691
$ws.Range("M2").Value = 27.937112
$ws.Range("N2").Value = 55.874224
$ws.Range("O2").Value = 0.1933230395895585
$ws.Range("P2").Value = 0.1911777832146586
$ws.Range("Q2").Value = 558.42731890498
$ws.Range("R2").Value = 2233.70927561992
$ws.Range("S2").Value = 0.05040416996064347
$ws.Range("T2").Value = 0.04177744391223413
$ws.Range("G3").Value = 19.9887275
$ws.Range("H3").Value = 39.977455
$ws.Range("I3").Value = 0.2607251058521317
$ws.Range("J3").Value = 0.2185266677421691
$ws.Range("O3").Value = 0.01881664736400372
$ws.Range("P3").Value = 0.02791176575399093
$ws.Range("Q3").Value = 54.35322122272584
$ws.Range("R3").Value = 326.119327336355
$ws.Range("S3").Value = 0.004905972375762103
$ws.Range("T3").Value = 0.006099465161019629
$ws.Range("G4").Value = 19.9887275
$ws.Range("H4").Value = 39.977455
$ws.Range("I4").Value = 0.2607251058521317
$ws.Range("J4").Value = 0.2185266677421691
$ws.Range("M4").Value = 0.08232966666666668
$ws.Range("N4").Value = 0.246989
$ws.Range("O4").Value = 0.0005697160611445862
$ws.Range("P4").Value = 0.0008450911013709168
$ws.Range("Q4").Value = 1.645665272165834
$ws.Range("R4").Value = 9.873991632995001
$ws.Range("S4").Value = 0.0001485392803475818
$ws.Range("T4").Value = 0.0001846749423211461
$ws.Range("G5").Value = 19.9887275
$ws.Range("H5").Value = 39.977455
$ws.Range("I5").Value = 0.2607251058521317
$ws.Range("J5").Value = 0.2185266677421691
$ws.Range("M5").Value = 113.329716
$ws.Range("N5").Value = 226.659432
$ws.Range("O5").Value = 0.7842344324259937
$ws.Range("P5").Value = 0.775531983306894
$ws.Range("Q5").Value = 2265.31681077639
$ws.Range("R5").Value = 9061.26724310556
$ws.Range("S5").Value = 0.2044696054071536
$ws.Range("T5").Value = 0.169474420039531
$ws.Range("G6").Value = 19.9887275
$ws.Range("H6").Value = 39.977455
$ws.Range("I6").Value = 0.2607251058521317
$ws.Range("J6").Value = 0.2185266677421691
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3519716666666666
$ws.Range("N6").Value = 1.055915
$ws.Range("O6").Value = 0.002435621564942105
$ws.Range("P6").Value = 0.003612891142132126
$ws.Range("Q6").Value = 7.035465732720833
$ws.Range("R6").Value = 42.21279439632499
$ws.Range("S6").Value = 0.0006350276903352649
$ws.Range("T6").Value = 0.0007895130622053327
$ws.Range("G7").Value = 19.9887275
$ws.Range("H7").Value = 39.977455
$ws.Range("I7").Value = 0.2607251058521317
$ws.Range("J7").Value = 0.2185266677421691
$ws.Range("M7").Value = 0.08967466666666667
$ws.Range("N7").Value = 0.269024
$ws.Range("O7").Value = 0.0006205429943574862
$ws.Range("P7").Value = 0.0009204854809534413
$ws.Range("Q7").Value = 1.792482475653333
$ws.Range("R7").Value = 10.75489485392
$ws.Range("S7").Value = 0.0001617911378896543
$ws.Range("T7").Value = 0.0002011506248578034
$ws.Range("I8").Value = 0.09121576982520653
$ws.Range("J8").Value = 0.1146787043788409
$ws.Range("M8").Value = 27.937112
$ws.Range("N8").Value = 55.874224
$ws.Range("O8").Value = 0.1933230395895585
$ws.Range("P8").Value = 0.1911777832146586
$ws.Range("Q8").Value = 195.36813541168
$ws.Range("R8").Value = 1172.20881247008
$ws.Range("S8").Value = 0.01763410988111045
$ws.Range("T8").Value = 0.02192402048507597
$ws.Range("I9").Value = 0.09121576982520653
$ws.Range("J9").Value = 0.1146787043788409
$ws.Range("O9").Value = 0.01881664736400372
$ws.Range("P9").Value = 0.02791176575399093
$ws.Range("S9").Value = 0.001716374974837042
$ws.Range("T9").Value = 0.003200885133593382
$ws.Range("I10").Value = 0.09121576982520653
$ws.Range("J10").Value = 0.1146787043788409
$ws.Range("M10").Value = 0.08232966666666668
$ws.Range("N10").Value = 0.246989
$ws.Range("O10").Value = 0.0005697160611445862
$ws.Range("P10").Value = 0.0008450911013709168
$ws.Range("Q10").Value = 0.5757428851533334
$ws.Range("R10").Value = 5.18168596638
$ws.Range("S10").Value = 0.00005196708909908786
$ws.Range("T10").Value = 0.00009691395258730447
$ws.Range("I11").Value = 0.09121576982520653
$ws.Range("J11").Value = 0.1146787043788409
$ws.Range("M11").Value = 113.329716
$ws.Range("N11").Value = 226.659432
$ws.Range("O11").Value = 0.7842344324259937
$ws.Range("P11").Value = 0.775531983306894
$ws.Range("Q11").Value = 792.5305701482399
$ws.Range("R11").Value = 4755.18342088944
$ws.Range("S11").Value = 0.07153454747717092
$ws.Range("T11").Value = 0.0889370030499875
$ws.Range("I12").Value = 0.09121576982520653
$ws.Range("J12").Value = 0.1146787043788409
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.3519716666666666
$ws.Range("N12").Value = 1.055915
$ws.Range("O12").Value = 0.002435621564942105
$ws.Range("P12").Value = 0.003612891142132126
$ws.Range("Q12").Value = 2.461387141033333
$ws.Range("R12").Value = 22.1524842693
$ws.Range("S12").Value = 0.0002221670960490684
$ws.Range("T12").Value = 0.000414321675241503
$ws.Range("I13").Value = 0.09121576982520653
$ws.Range("J13").Value = 0.1146787043788409
$ws.Range("M13").Value = 0.08967466666666667
$ws.Range("N13").Value = 0.269024
$ws.Range("O13").Value = 0.0006205429943574862
$ws.Range("P13").Value = 0.0009204854809534413
$ws.Range("Q13").Value = 0.6271074984533332
$ws.Range("R13").Value = 5.643967486079999
$ws.Range("S13").Value = 0.0000566033069399569
$ws.Range("T13").Value = 0.0001055600823552749
$ws.Range("G14").Value = 4.517918
$ws.Range("H14").Value = 13.553754
$ws.Range("I14").Value = 0.0589299468303448
$ws.Range("J14").Value = 0.07408817537327214
$ws.Range("M14").Value = 27.937112
$ws.Range("N14").Value = 55.874224
$ws.Range("O14").Value = 0.1933230395895585
$ws.Range("P14").Value = 0.1911777832146586
$ws.Range("Q14").Value = 126.217581172816
$ws.Range("R14").Value = 757.3054870368959
$ws.Range("S14").Value = 0.01139251644409332
$ws.Range("T14").Value = 0.01416401313028103
$ws.Range("G15").Value = 4.517918
$ws.Range("H15").Value = 13.553754
$ws.Range("I15").Value = 0.0589299468303448
$ws.Range("J15").Value = 0.07408817537327214
$ws.Range("O15").Value = 0.01881664736400372
$ws.Range("P15").Value = 0.02791176575399093
$ws.Range("Q15").Value = 12.28509401211933
$ws.Range("R15").Value = 110.565846109074
$ws.Range("S15").Value = 0.001108864028686087
$ws.Range("T15").Value = 0.002067931796159371
$ws.Range("G16").Value = 4.517918
$ws.Range("H16").Value = 13.553754
$ws.Range("I16").Value = 0.0589299468303448
$ws.Range("J16").Value = 0.07408817537327214
$ws.Range("M16").Value = 0.08232966666666668
$ws.Range("N16").Value = 0.246989
$ws.Range("O16").Value = 0.0005697160611445862
$ws.Range("P16").Value = 0.0008450911013709168
$ws.Range("Q16").Value = 0.3719586829673334
$ws.Range("R16").Value = 3.347628146706
$ws.Range("S16").Value = 0.00003357333719164394
$ws.Range("T16").Value = 0.00006261125772476019
$ws.Range("G17").Value = 4.517918
$ws.Range("H17").Value = 13.553754
$ws.Range("I17").Value = 0.0589299468303448
$ws.Range("J17").Value = 0.07408817537327214
$ws.Range("M17").Value = 113.329716
$ws.Range("N17").Value = 226.659432
$ws.Range("O17").Value = 0.7842344324259937
$ws.Range("P17").Value = 0.775531983306894
$ws.Range("Q17").Value = 512.014363851288
$ws.Range("R17").Value = 3072.086183107728
$ws.Range("S17").Value = 0.04621489340538944
$ws.Range("T17").Value = 0.05745774958682272
$ws.Range("G18").Value = 4.517918
$ws.Range("H18").Value = 13.553754
$ws.Range("I18").Value = 0.0589299468303448
$ws.Range("J18").Value = 0.07408817537327214
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.3519716666666666
$ws.Range("N18").Value = 1.055915
$ws.Range("O18").Value = 0.002435621564942105
$ws.Range("P18").Value = 0.003612891142132126
$ws.Range("Q18").Value = 1.590179128323333
$ws.Range("R18").Value = 14.31161215491
$ws.Range("S18").Value = 0.0001435310493208795
$ws.Range("T18").Value = 0.0002676725125428264
$ws.Range("G19").Value = 4.517918
$ws.Range("H19").Value = 13.553754
$ws.Range("I19").Value = 0.0589299468303448
$ws.Range("J19").Value = 0.07408817537327214
$ws.Range("M19").Value = 0.08967466666666667
$ws.Range("N19").Value = 0.269024
$ws.Range("O19").Value = 0.0006205429943574862
$ws.Range("P19").Value = 0.0009204854809534413
$ws.Range("Q19").Value = 0.4051427906773333
$ws.Range("R19").Value = 3.646285116096
$ws.Range("S19").Value = 0.00003656856566342962
$ws.Range("T19").Value = 0.00006819708974142932
$ws.Range("G20").Value = 27.0681495
$ws.Range("H20").Value = 54.136299
$ws.Range("I20").Value = 0.3530663041761326
$ws.Range("J20").Value = 0.2959224148801799
$ws.Range("M20").Value = 27.937112
$ws.Range("N20").Value = 55.874224
$ws.Range("O20").Value = 0.1933230395895585
$ws.Range("P20").Value = 0.1911777832146586
$ws.Range("Q20").Value = 756.205924214244
$ws.Range("R20").Value = 3024.823696856976
$ws.Range("S20").Value = 0.06825585109998156
$ws.Range("T20").Value = 0.05657379128032129
$ws.Range("G21").Value = 27.0681495
$ws.Range("H21").Value = 54.136299
$ws.Range("I21").Value = 0.3530663041761326
$ws.Range("J21").Value = 0.2959224148801799
$ws.Range("O21").Value = 0.01881664736400372
$ws.Range("P21").Value = 0.02791176575399093
$ws.Range("Q21").Value = 73.60354068878651
$ws.Range("R21").Value = 441.621244132719
$ws.Range("S21").Value = 0.00664352414179436
$ws.Range("T21").Value = 0.0082597171254909
$ws.Range("G22").Value = 27.0681495
$ws.Range("H22").Value = 54.136299
$ws.Range("I22").Value = 0.3530663041761326
$ws.Range("J22").Value = 0.2959224148801799
$ws.Range("M22").Value = 0.08232966666666668
$ws.Range("N22").Value = 0.246989
$ws.Range("O22").Value = 0.0005697160611445862
$ws.Range("P22").Value = 0.0008450911013709168
$ws.Range("Q22").Value = 2.2285117256185
$ws.Range("R22").Value = 13.371070353711
$ws.Range("S22").Value = 0.0002011475441381026
$ws.Range("T22").Value = 0.0002500813995114326
$ws.Range("G23").Value = 27.0681495
$ws.Range("H23").Value = 54.136299
$ws.Range("I23").Value = 0.3530663041761326
$ws.Range("J23").Value = 0.2959224148801799
$ws.Range("M23").Value = 113.329716
$ws.Range("N23").Value = 226.659432
$ws.Range("O23").Value = 0.7842344324259937
$ws.Range("P23").Value = 0.775531983306894
$ws.Range("Q23").Value = 3067.625695480542
$ws.Range("R23").Value = 12270.50278192217
$ws.Range("S23").Value = 0.2768867526643126
$ws.Range("T23").Value = 0.2294972973169914
$ws.Range("G24").Value = 27.0681495
$ws.Range("H24").Value = 54.136299
$ws.Range("I24").Value = 0.3530663041761326
$ws.Range("J24").Value = 0.2959224148801799
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.3519716666666666
$ws.Range("N24").Value = 1.055915
$ws.Range("O24").Value = 0.002435621564942105
$ws.Range("P24").Value = 0.003612891142132126
$ws.Range("Q24").Value = 9.527221693097498
$ws.Range("R24").Value = 57.163330158585
$ws.Range("S24").Value = 0.0008599359043057973
$ws.Range("T24").Value = 0.00106913547147895
$ws.Range("G25").Value = 27.0681495
$ws.Range("H25").Value = 54.136299
$ws.Range("I25").Value = 0.3530663041761326
$ws.Range("J25").Value = 0.2959224148801799
$ws.Range("M25").Value = 0.08967466666666667
$ws.Range("N25").Value = 0.269024
$ws.Range("O25").Value = 0.0006205429943574862
$ws.Range("P25").Value = 0.0009204854809534413
$ws.Range("Q25").Value = 2.427327283696
$ws.Range("R25").Value = 14.563963702176
$ws.Range("S25").Value = 0.0002190928216001883
$ws.Range("T25").Value = 0.0002723922863858862
$ws.Range("G26").Value = 5.535821333333334
$ws.Range("H26").Value = 16.607464
$ws.Range("I26").Value = 0.0722070778698555
$ws.Range("J26").Value = 0.09078051035434934
$ws.Range("M26").Value = 27.937112
$ws.Range("N26").Value = 55.874224
$ws.Range("O26").Value = 0.1933230395895585
$ws.Range("P26").Value = 0.1911777832146586
$ws.Range("Q26").Value = 154.6548606013227
$ws.Range("R26").Value = 927.9291636079359
$ws.Range("S26").Value = 0.0139592917736804
$ws.Range("T26").Value = 0.01735521672863987
$ws.Range("G27").Value = 5.535821333333334
$ws.Range("H27").Value = 16.607464
$ws.Range("I27").Value = 0.0722070778698555
$ws.Range("J27").Value = 0.09078051035434934
$ws.Range("O27").Value = 0.01881664736400372
$ws.Range("P27").Value = 0.02791176575399093
$ws.Range("Q27").Value = 15.05297030939822
$ws.Range("R27").Value = 135.476732784584
$ws.Range("S27").Value = 0.001358695121462228
$ws.Range("T27").Value = 0.002533844340038346
$ws.Range("G28").Value = 5.535821333333334
$ws.Range("H28").Value = 16.607464
$ws.Range("I28").Value = 0.0722070778698555
$ws.Range("J28").Value = 0.09078051035434934
$ws.Range("M28").Value = 0.08232966666666668
$ws.Range("N28").Value = 0.246989
$ws.Range("O28").Value = 0.0005697160611445862
$ws.Range("P28").Value = 0.0008450911013709168
$ws.Range("Q28").Value = 0.4557623250995557
$ws.Range("R28").Value = 4.101860925896
$ws.Range("S28").Value = 0.00004113753199077449
$ws.Range("T28").Value = 0.00007671780147837099
$ws.Range("G29").Value = 5.535821333333334
$ws.Range("H29").Value = 16.607464
$ws.Range("I29").Value = 0.0722070778698555
$ws.Range("J29").Value = 0.09078051035434934
$ws.Range("M29").Value = 113.329716
$ws.Range("N29").Value = 226.659432
$ws.Range("O29").Value = 0.7842344324259937
$ws.Range("P29").Value = 0.775531983306894
$ws.Range("Q29").Value = 627.3730595334081
$ws.Range("R29").Value = 3764.238357200448
$ws.Range("S29").Value = 0.05662727673040566
$ws.Range("T29").Value = 0.07040318924072057
$ws.Range("G30").Value = 5.535821333333334
$ws.Range("H30").Value = 16.607464
$ws.Range("I30").Value = 0.0722070778698555
$ws.Range("J30").Value = 0.09078051035434934
$ws.Range("K30").Value = 2
$ws.Range("L30").Value = 0.6666666666666666
$ws.Range("M30").Value = 0.3519716666666666
$ws.Range("N30").Value = 1.055915
$ws.Range("O30").Value = 0.002435621564942105
$ws.Range("P30").Value = 0.003612891142132126
$ws.Range("Q30").Value = 1.948452261062222
$ws.Range("R30").Value = 17.53607034956
$ws.Range("S30").Value = 0.0001758691160012739
$ws.Range("T30").Value = 0.0003279801017374624
$ws.Range("G31").Value = 5.535821333333334
$ws.Range("H31").Value = 16.607464
$ws.Range("I31").Value = 0.0722070778698555
$ws.Range("J31").Value = 0.09078051035434934
$ws.Range("M31").Value = 0.08967466666666667
$ws.Range("N31").Value = 0.269024
$ws.Range("O31").Value = 0.0006205429943574862
$ws.Range("P31").Value = 0.0009204854809534413
$ws.Range("Q31").Value = 0.4964229327928889
$ws.Range("R31").Value = 4.467806395136
$ws.Range("S31").Value = 0.00004480759631516431
$ws.Range("T31").Value = 0.00008356214173472211
$ws.Range("G32").Value = 12.56215366666667
$ws.Range("H32").Value = 37.68646099999999
$ws.Range("I32").Value = 0.163855795446329
$ws.Range("J32").Value = 0.2060035272711885
$ws.Range("M32").Value = 27.937112
$ws.Range("N32").Value = 55.874224
$ws.Range("O32").Value = 0.1933230395895585
$ws.Range("P32").Value = 0.1911777832146586
$ws.Range("Q32").Value = 350.9502939468773
$ws.Range("R32").Value = 2105.701763681263
$ws.Range("S32").Value = 0.03167710043004925
$ws.Range("T32").Value = 0.03938329767810629
$ws.Range("G33").Value = 12.56215366666667
$ws.Range("H33").Value = 37.68646099999999
$ws.Range("I33").Value = 0.163855795446329
$ws.Range("J33").Value = 0.2060035272711885
$ws.Range("O33").Value = 0.01881664736400372
$ws.Range("P33").Value = 0.02791176575399093
$ws.Range("Q33").Value = 34.15892869009345
$ws.Range("R33").Value = 307.430358210841
$ws.Range("S33").Value = 0.003083216721461899
$ws.Range("T33").Value = 0.005749922197689297
$ws.Range("G34").Value = 12.56215366666667
$ws.Range("H34").Value = 37.68646099999999
$ws.Range("I34").Value = 0.163855795446329
$ws.Range("J34").Value = 0.2060035272711885
$ws.Range("M34").Value = 0.08232966666666668
$ws.Range("N34").Value = 0.246989
$ws.Range("O34").Value = 0.0005697160611445862
$ws.Range("P34").Value = 0.0008450911013709168
$ws.Range("Q34").Value = 1.034237923992111
$ws.Range("R34").Value = 9.308141315928999
$ws.Range("S34").Value = 0.00009335127837739557
$ws.Range("T34").Value = 0.0001740917477479024
$ws.Range("G35").Value = 12.56215366666667
$ws.Range("H35").Value = 37.68646099999999
$ws.Range("I35").Value = 0.163855795446329
$ws.Range("J35").Value = 0.2060035272711885
$ws.Range("M35").Value = 113.329716
$ws.Range("N35").Value = 226.659432
$ws.Range("O35").Value = 0.7842344324259937
$ws.Range("P35").Value = 0.775531983306894
$ws.Range("Q35").Value = 1423.665307391692
$ws.Range("R35").Value = 8541.991844350152
$ws.Range("S35").Value = 0.1285013567415615
$ws.Range("T35").Value = 0.1597623240728407
$ws.Range("G36").Value = 12.56215366666667
$ws.Range("H36").Value = 37.68646099999999
$ws.Range("I36").Value = 0.163855795446329
$ws.Range("J36").Value = 0.2060035272711885
$ws.Range("K36").Value = 2
$ws.Range("L36").Value = 0.6666666666666666
$ws.Range("M36").Value = 0.3519716666666666
$ws.Range("N36").Value = 1.055915
$ws.Range("O36").Value = 0.002435621564942105
$ws.Range("P36").Value = 0.003612891142132126
$ws.Range("Q36").Value = 4.421522162979444
$ws.Range("R36").Value = 39.79369946681499
$ws.Range("S36").Value = 0.0003990907089298212
$ws.Range("T36").Value = 0.0007442683189260508
$ws.Range("G37").Value = 12.56215366666667
$ws.Range("H37").Value = 37.68646099999999
$ws.Range("I37").Value = 0.163855795446329
$ws.Range("J37").Value = 0.2060035272711885
$ws.Range("M37").Value = 0.08967466666666667
$ws.Range("N37").Value = 0.269024
$ws.Range("O37").Value = 0.0006205429943574862
$ws.Range("P37").Value = 0.0009204854809534413
$ws.Range("Q37").Value = 1.126506942673778
$ws.Range("R37").Value = 10.138562484064
$ws.Range("S37").Value = 0.0001016795659490927
$ws.Range("T37").Value = 0.0001896232558783253
